$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

$wsALC.Range("H15").Value = 1140.662
$wsALC.Range("I15").Value = 1140.662
$wsALC.Range("K15").Value = 3421.986
$wsALC.Range("M15").Value = -3252.986
$wsALC.Range("H53").Value = 37621.035
$wsALC.Range("J53").Value = 71714.28999999999
$wsALC.Range("L53").Value = 71714.28999999999
$wsALC.Range("N53").Value = -72988.28999999999
$wsALC.Range("H55").Value = 268.66666
$wsALC.Range("I55").Value = 236.66667
$wsALC.Range("J55").Value = 300.66666
$wsALC.Range("K55").Value = 236.66667
$wsALC.Range("L55").Value = 300.66666
$wsALC.Range("M55").Value = -22.66667000000001
$wsALC.Range("N55").Value = -728.66666
$wsALC.Range("H137").Value = 4637.522
$wsALC.Range("I137").Value = 3674.2307
$wsALC.Range("K137").Value = 11022.6921
$wsALC.Range("M137").Value = -8472.6921
$wsARM.Range("H2").Value = 628.73914
$wsARM.Range("I2").Value = 612.4211
$wsARM.Range("K2").Value = 612.4211
$wsARM.Range("M2").Value = -499.4211
$wsARM.Range("H5").Value = 1070
$wsARM.Range("I5").Value = 341.25
$wsARM.Range("K5").Value = 341.25
$wsARM.Range("M5").Value = -229.25
$wsARM.Range("H28").Value = 18671.7
$wsARM.Range("I28").Value = 17968.555
$wsARM.Range("K28").Value = 17968.555
$wsARM.Range("M28").Value = -17776.555
$wsARM.Range("H32").Value = 10791.432
$wsARM.Range("I32").Value = 6673.2686
$wsARM.Range("K32").Value = 6673.2686
$wsARM.Range("M32").Value = -6386.2686
$wsARM.Range("H45").Value = 506539.47
$wsARM.Range("I45").Value = 695788.6
$wsARM.Range("K45").Value = 695788.6
$wsARM.Range("M45").Value = -695411.6
$wsARM.Range("H74").Value = 41670076
$wsARM.Range("I74").Value = 55558030
$wsARM.Range("K74").Value = 55558030
$wsARM.Range("M74").Value = -55557156
$wsARM.Range("H77").Value = 41670076
$wsARM.Range("I77").Value = 55558030
$wsARM.Range("K77").Value = 277790150
$wsARM.Range("M77").Value = -277785782
$wsARM.Range("H99").Value = 18671.7
$wsARM.Range("I99").Value = 17968.555
$wsARM.Range("K99").Value = 17968.555
$wsARM.Range("M99").Value = -14973.555
$wsARM.Range("H116").Value = 628.73914
$wsARM.Range("I116").Value = 612.4211
$wsARM.Range("K116").Value = 612.4211
$wsARM.Range("M116").Value = 1681.5789
$wsARM.Range("H122").Value = 5610.625
$wsARM.Range("I122").Value = 4962.3335
$wsARM.Range("K122").Value = 14887.0005
$wsARM.Range("M122").Value = -12437.0005
$wsBSM.Range("H3").Value = 628.73914
$wsBSM.Range("I3").Value = 612.4211
$wsBSM.Range("K3").Value = 612.4211
$wsBSM.Range("M3").Value = -498.4211
$wsBSM.Range("H4").Value = 1070
$wsBSM.Range("I4").Value = 341.25
$wsBSM.Range("K4").Value = 341.25
$wsBSM.Range("M4").Value = -226.25
$wsBSM.Range("H94").Value = 6251109
$wsBSM.Range("I94").Value = 6945513
$wsBSM.Range("K94").Value = 6945513
$wsBSM.Range("M94").Value = -6945062
$wsBSM.Range("H99").Value = 3134.7058
$wsBSM.Range("I99").Value = 2537.0833
$wsBSM.Range("J99").Value = 4569
$wsBSM.Range("K99").Value = 2537.0833
$wsBSM.Range("L99").Value = 4569
$wsBSM.Range("M99").Value = -1039.0833
$wsBSM.Range("N99").Value = -7565
$wsCRP.Range("H41").Value = 43314.25
$wsCRP.Range("I41").Value = 24259
$wsCRP.Range("J41").Value = 49666
$wsCRP.Range("K41").Value = 24259
$wsCRP.Range("L41").Value = 49666
$wsCRP.Range("M41").Value = -23831
$wsCRP.Range("N41").Value = -50522
$wsCRP.Range("H62").Value = 103437
$wsCRP.Range("I62").Value = 3869.375
$wsCRP.Range("J62").Value = 217228.58
$wsCRP.Range("K62").Value = 3869.375
$wsCRP.Range("L62").Value = 217228.58
$wsCRP.Range("M62").Value = -3245.375
$wsCRP.Range("N62").Value = -218476.58
$wsCRP.Range("H65").Value = 103437
$wsCRP.Range("I65").Value = 3869.375
$wsCRP.Range("J65").Value = 217228.58
$wsCRP.Range("K65").Value = 19346.875
$wsCRP.Range("L65").Value = 1086142.9
$wsCRP.Range("M65").Value = -16226.875
$wsCRP.Range("N65").Value = -1092382.9
$wsCUL.Range("H75").Value = 2900
$wsCUL.Range("J75").Value = 4250
$wsCUL.Range("L75").Value = 12750
$wsCUL.Range("N75").Value = -14746
$wsCUL.Range("H78").Value = 2900
$wsCUL.Range("J78").Value = 4250
$wsCUL.Range("L78").Value = 38250
$wsCUL.Range("N78").Value = -48234
$wsCUL.Range("H80").Value = 5613.636
$wsCUL.Range("J80").Value = 6000
$wsCUL.Range("L80").Value = 18000
$wsCUL.Range("N80").Value = -19872
$wsCUL.Range("H83").Value = 5613.636
$wsCUL.Range("J83").Value = 6000
$wsCUL.Range("L83").Value = 54000
$wsCUL.Range("N83").Value = -63360
$wsCUL.Range("H92").Value = 749.25
$wsCUL.Range("I92").Value = 499.5
$wsCUL.Range("K92").Value = 1498.5
$wsCUL.Range("M92").Value = -250.5
$wsCUL.Range("H98").Value = 920.4
$wsCUL.Range("J98").Value = 920.4
$wsCUL.Range("L98").Value = 2761.2
$wsCUL.Range("N98").Value = -5757.2
$wsCUL.Range("H103").Value = 1062.5
$wsCUL.Range("I103").Value = 743.75
$wsCUL.Range("K103").Value = 2231.25
$wsCUL.Range("M103").Value = -1352.25
$wsCUL.Range("H139").Value = 1672896.2
$wsCUL.Range("I139").Value = 1967725.1
$wsCUL.Range("K139").Value = 5903175.300000001
$wsCUL.Range("M139").Value = -5898035.300000001
$wsGSM.Range("H80").Value = 11068.529
$wsGSM.Range("J80").Value = 13781.923
$wsGSM.Range("L80").Value = 13781.923
$wsGSM.Range("N80").Value = -15777.923
$wsGSM.Range("H83").Value = 11068.529
$wsGSM.Range("J83").Value = 13781.923
$wsGSM.Range("L83").Value = 68909.61500000001
$wsGSM.Range("N83").Value = -78893.61500000001
$wsGSM.Range("H107").Value = 471.5625
$wsGSM.Range("I107").Value = 443
$wsGSM.Range("J107").Value = 900
$wsGSM.Range("K107").Value = 443
$wsGSM.Range("L107").Value = 900
$wsGSM.Range("M107").Value = 1477
$wsGSM.Range("N107").Value = -4740
$wsLTW.Range("H7").Value = 25004262
$wsLTW.Range("I7").Value = 83336500
$wsLTW.Range("J7").Value = 4732.7144
$wsLTW.Range("K7").Value = 83336500
$wsLTW.Range("L7").Value = 4732.7144
$wsLTW.Range("M7").Value = -83336388
$wsLTW.Range("N7").Value = -4956.7144
$wsLTW.Range("H46").Value = 2387.3333
$wsLTW.Range("I46").Value = 1639.4
$wsLTW.Range("J46").Value = 2761.3
$wsLTW.Range("K46").Value = 1639.4
$wsLTW.Range("L46").Value = 2761.3
$wsLTW.Range("M46").Value = -1451.4
$wsLTW.Range("N46").Value = -3137.3
$wsLTW.Range("H55").Value = 969.8570999999999
$wsLTW.Range("I55").Value = 494
$wsLTW.Range("K55").Value = 494
$wsLTW.Range("M55").Value = -321
$wsLTW.Range("H100").Value = 5553.8887
$wsLTW.Range("I100").Value = 4995
$wsLTW.Range("K100").Value = 4995
$wsLTW.Range("M100").Value = -4454
$wsLTW.Range("H109").Value = 150001
$wsLTW.Range("J109").Value = 150001
$wsLTW.Range("L109").Value = 150001
$wsLTW.Range("N109").Value = -152775
$wsLTW.Range("H126").Value = 25004262
$wsLTW.Range("I126").Value = 83336500
$wsLTW.Range("J126").Value = 4732.7144
$wsLTW.Range("K126").Value = 250009500
$wsLTW.Range("L126").Value = 14198.1432
$wsLTW.Range("M126").Value = -250007030
$wsLTW.Range("N126").Value = -19138.1432
$wsLTW.Range("H132").Value = 17246290
$wsLTW.Range("I132").Value = 27780908
$wsLTW.Range("J132").Value = 7826
$wsLTW.Range("K132").Value = 83342724
$wsLTW.Range("L132").Value = 23478
$wsLTW.Range("M132").Value = -83340194
$wsLTW.Range("N132").Value = -28538
$wsLTW.Range("H136").Value = 6492.558
$wsLTW.Range("I136").Value = 6637.4053
$wsLTW.Range("J136").Value = 5599.3335
$wsLTW.Range("K136").Value = 19912.2159
$wsLTW.Range("L136").Value = 16798.0005
$wsLTW.Range("M136").Value = -17362.2159
$wsLTW.Range("N136").Value = -21898.0005
$wsWVR.Range("H11").Value = 3150
$wsWVR.Range("J11").Value = 4000
$wsWVR.Range("L11").Value = 4000
$wsWVR.Range("N11").Value = -4284
$wsWVR.Range("H39").Value = 0
$wsWVR.Range("I39").Value = 0
$wsWVR.Range("K39").Value = 0
$wsWVR.Range("M39").ClearContents()
$wsWVR.Range("H61").Value = 7000
$wsWVR.Range("I61").Value = 7000
$wsWVR.Range("K61").Value = 7000
$wsWVR.Range("M61").Value = -6708
$wsWVR.Range("H96").Value = 2644.3333
$wsWVR.Range("I96").Value = 1623.3334
$wsWVR.Range("J96").Value = 3665.3333
$wsWVR.Range("K96").Value = 1623.3334
$wsWVR.Range("L96").Value = 3665.3333
$wsWVR.Range("M96").Value = -250.3334
$wsWVR.Range("N96").Value = -6411.3333
$wsWVR.Range("H132").Value = 15377.064
$wsWVR.Range("I132").Value = 5581.2593
$wsWVR.Range("K132").Value = 16743.7779
$wsWVR.Range("M132").Value = -14213.7779
$wsWVR.Range("H136").Value = 5365.763
$wsWVR.Range("I136").Value = 8060
$wsWVR.Range("K136").Value = 24180
$wsWVR.Range("M136").Value = -21630
